# Generate Report for Handback
# For the source file "107cef57-d204-4b90-a9f7-881e77a74b86.md" (row 6 of the
# status tables), the handback transform failed because the handback type
# (mt) did not match the handoff type (ht). Reflect this in the
# Overview/zh-cn/de-de tables: update the Status text and record the new
# Error Detail message for each localized sheet.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handback transform failed"
$errorDetail = "The handback type mt is not match with handoff type ht."

# Overview sheet: zh-cn and de-de status columns (B6/C6) for this file.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B6").Value = $newStatus
$overview.Range("C6").Value = $newStatus

# zh-cn sheet: update Status (C6) and record the Error Detail (L6).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = $newStatus
$zhcn.Range("L6").Value = $errorDetail

# de-de sheet: update Status (C6) and record the Error Detail (L6).
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = $newStatus
$dede.Range("L6").Value = $errorDetail
